# ---------------------------------------------------------------------------
# Add envelopes to the loadcase generation script
#   * Rename Sheet1 -> Loadcases
#   * Add a new worksheet "Envelopes" after Loadcases
#   * Add a "Gravity" column (C) to the Loadcases sheet, shifting the
#     existing "Analysis" column from C to D
#   * Add a "Prestress" loadcase row and notes column (I) to Loadcases
#   * Populate the Envelopes sheet with Settlement / Wind / Temperature
#     envelopes, the "Other" traffic loadcases 1-10 and notes column (G)
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- rename the existing sheet -------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Loadcases"

# --- add the new Envelopes sheet after Loadcases --------------------------
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Envelopes"

# ===========================================================================
# Loadcases sheet
# ===========================================================================

# Header row -----------------------------------------------------------------
$ws1.Range("A1").Value = "Name"
$ws1.Range("B1").Value = "Count"
$ws1.Range("C1").Value = "Gravity"
$ws1.Range("D1").Value = "Analysis"
$ws1.Range("A1").Copy()
$ws1.Range("D1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Move the "Analysis" note text from C8 to D8 (Gravity column inserted at C)
$ws1.Range("D8").Value = "Moving Load Characteristic"
$ws1.Range("C8").ClearContents()

# New Gravity column values
$ws1.Range("C2").Value = "Yes"
$ws1.Range("C11").Value = "Yes"

# New Prestress loadcase row
$ws1.Range("A11").Value = "Prestress"
$ws1.Range("D11").Value = "Prestress"

# Notes column (I) - red text
$ws1.Range("I2").Value = "Name = Name of loadcase to be created"
$ws1.Range("I3").Value = "Count = Number of similar loadcases to be created = default empty = 1"
$ws1.Range("I4").Value = "Gravity = Automatically apply gravity to the loadcase"
$ws1.Range("I5").Value = "Analysis = Name of the analysis in which to create the loadcase"
$ws1.Range("I2:I5").Font.Color = 255

# Centre-align the Gravity / Count columns (B:C)
$ws1.Range("B5").HorizontalAlignment = -4108
$ws1.Range("B7").HorizontalAlignment = -4108
$ws1.Range("B8").HorizontalAlignment = -4108
$ws1.Range("B9").HorizontalAlignment = -4108
$ws1.Range("C2").HorizontalAlignment = -4108
$ws1.Range("C11").HorizontalAlignment = -4108

# column widths
$ws1.Columns.Item(2).ColumnWidth = 8
$ws1.Columns.Item(3).ColumnWidth = 8

$ws1.Range("D22").Select() | Out-Null

# ===========================================================================
# Envelopes sheet
# ===========================================================================

$ws2.Range("A1").Value = "Name"
$ws2.Range("B1").Value = "Loadcases"
$ws2.Range("C1").Value = "FindSimilar"

$ws2.Range("A2").Value = "Settlement"
$ws2.Range("B2").Value = "Settlement"
$ws2.Range("C2").Value = "Yes"

$ws2.Range("A3").Value = "Wind"
$ws2.Range("B3").Value = "Wind"
$ws2.Range("C3").Value = "Yes"

$ws2.Range("A4").Value = "Temperature"
$ws2.Range("B4").Value = "Temperature"
$ws2.Range("C4").Value = "Yes"

$ws2.Range("A5").Value = "Other"
$ws2.Range("B5").Value = "Traffic 1"
$ws2.Range("A6").Value = "Other"
$ws2.Range("B6").Value = "Traffic 2"
$ws2.Range("A7").Value = "Other"
$ws2.Range("B7").Value = "Traffic 3"
$ws2.Range("A8").Value = "Other"
$ws2.Range("B8").Value = "Traffic 4"
$ws2.Range("A9").Value = "Other"
$ws2.Range("B9").Value = "Traffic 5"
$ws2.Range("A10").Value = "Other"
$ws2.Range("B10").Value = "Traffic 6"
$ws2.Range("A11").Value = "Other"
$ws2.Range("B11").Value = "Traffic 7"
$ws2.Range("A12").Value = "Other"
$ws2.Range("B12").Value = "Traffic 8"
$ws2.Range("A13").Value = "Other"
$ws2.Range("B13").Value = "Traffic 9"
$ws2.Range("A14").Value = "Other"
$ws2.Range("B14").Value = "Traffic 10"

# Notes column (G) - red text
$ws2.Range("G2").Value = "Name = Name of envelope to be created"
$ws2.Range("G3").Value = "Loadcases = Name of loadcases to be included in the envelope"
$ws2.Range("G4").Value = "FindSimilar = Use all loadsets that start with the loadcase name"
$ws2.Range("G2:G4").Font.Color = 255

# Column widths (closest achievable given engine's 1/6-character rounding)
$ws2.Columns.Item(1).ColumnWidth = 14.833333333333334
$ws2.Columns.Item(2).ColumnWidth = 16.333333333333332
$ws2.Columns.Item(3).ColumnWidth = 21.833333333333332

$ws2.Range("G2:G4").Select() | Out-Null

# --- re-activate the Loadcases sheet so tab selection / active sheet match -
$ws1.Activate()
$ws1.Range("D22").Select() | Out-Null
